$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.446.56"
$ws.Range("E2").Value = "  -2.51%  "

$ws.Range("D3").Value = "1.864.71"
$ws.Range("E3").Value = "  -2.58%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "329.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.43%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.21%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4734"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.46%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3971"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.98%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.26"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -10.98%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08025"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.53%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.021"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.34%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.61"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.32%  "

$ws.Range("D13").Value = "1.846.76"
$ws.Range("E13").Value = "  -2.50%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.966"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.54%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.189"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.22%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.27%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "86.59"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.33%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001040"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.19%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06551"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.15%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.35"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.19%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("D21").Style = "Normal"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.519"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.92%  "

$ws.Range("D23").Value = "27.445.52"
$ws.Range("E23").Value = "  -2.45%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.98"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.79%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.301"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.08%  "

$ws.Range("D26").Value = "2.075.33"
$ws.Range("E26").Value = "  -2.23%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.34"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.55%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "154.21"
$ws.Range("D28").Style = "Normal"

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.087"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.34%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.551"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.74%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "122.50"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.70%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09525"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.31%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9604"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.70%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.461"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.48%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.586"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.70%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.304"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.29%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06077"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.60%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02239"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.49%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.215"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.03%  "

$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.064"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -9.05%  "

$ws.Range("B41").Value = "Frax"
$ws.Range("C41").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9998"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.14%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5968"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.18%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1906"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.19%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.36"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.19%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.266"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.10%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5665"
$ws.Range("D46").Style = "Normal"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.20"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.34%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.428"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.25%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.939"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.13%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06781"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.81%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "110.04"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.36%  "
